$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "CMS Username"
$ws.Range("B6").Value = "CMS Password"
$ws.Range("A7").Value = "Vamsi"
$ws.Range("B7").Value = "Vamsi123*"

$ws.Rows.Item(7).Select()
